{"js": "// The document originally stored the Title, Author and Abstract\n// paragraphs as a sequence of many runs (one run per word, plus a\n// separate run for every single space). The edit simply collapses each\n// of those paragraphs down to a single run holding the full text,\n// without changing the visible text or paragraph styles.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Styles of the three paragraphs that need their runs merged.\nconst stylesToMerge = [\"Title\", \"Author\", \"Abstract\"];\n\nfor (const paragraph of paragraphs.items) {\n  if (stylesToMerge.indexOf(paragraph.style) !== -1) {\n    // Re-inserting the paragraph's own text and replacing its current\n    // contents collapses all of the little single-word/space runs into\n    // one run containing the whole string.\n    paragraph.insertText(paragraph.text, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The Title, Author and Abstract paragraphs were originally split into\n# many runs (one run per word, plus a separate run for every single\n# space). Collapse each of those paragraphs down to a single run that\n# contains the whole string, without touching the visible text or the\n# paragraph styles.\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $styleName = $p.Range.Style.NameLocal\n    if ($styleName -eq \"Title\" -or $styleName -eq \"Author\" -or $styleName -eq \"Abstract\") {\n        $r = $p.Range\n        $r.MoveEnd(1, -1)   # exclude the paragraph mark from the range\n        $text = $r.Text\n        # Clearing the range before re-inserting the text forces Word to\n        # rebuild the paragraph with a single run instead of keeping the\n        # pre-existing run boundaries.\n        $r.Text = \"\"\n        $r.Text = $text\n    }\n}\n"}
